$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formats from row 2 onto the styled columns (E:H) in
# row 3 first, so the date/time/percent formatting carries over, same
# as the original authored workbook. Do this before entering the
# formulas below so the formula write isn't clobbered afterwards.
$ws.Range("E2:H2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3: mirror row 2 values via IF(ISBLANK(...),"",...) formulas.
# Column A gets its own (non-shared) formula.
$ws.Range("A3").Formula = '=IF(ISBLANK(A2),"",A2)'

# Columns B through I share the same formula pattern (relative refs
# adjust per-cell, same as typing the formula once and filling right).
$ws.Range("B3:I3").Formula = '=IF(ISBLANK(B2),"",B2)'

# Restore the active selection to C3, as it ends up after entering
# the new row of formulas.
$ws.Range("C3").Select()
